# Trade #10 closed at 2026-02-16 21:22:12 - leadlag DOWN +0.000%
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$ws.Cells.Item(9, 1).Value = 10

# Date-looking text must not be auto-converted to a date serial.
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "2026-02-16"
$ws.Cells.Item(9, 2).ClearFormats()

$ws.Cells.Item(9, 3).Value = "21:22:12"
$ws.Cells.Item(9, 4).Value = "leadlag"
$ws.Cells.Item(9, 5).Value = "DOWN"
$ws.Cells.Item(9, 6).Value = 69319.3
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 8).Value = "OPEN"
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0.7119
$ws.Cells.Item(9, 12).Value = "Binance leading with -0.071% move"
$ws.Cells.Item(9, 13).Value = ""
$ws.Cells.Item(9, 14).Value = 0
